# Switch the active/selected sheet from "infrastructure" to "o&m" and
# highlight row 26 on "o&m" with a solid red fill (the row being reviewed),
# updating the selections/scroll position on both sheets to match where the
# author was last looking ("update prints to follow what is going on").

$wb  = $excel.ActiveWorkbook
$omSheet    = $wb.Worksheets.Item("o&m")
$infraSheet = $wb.Worksheets.Item("infrastructure")

# Update "infrastructure"'s remembered selection first (selecting a range
# implicitly activates its sheet, so do this before activating "o&m" so the
# workbook ends up with "o&m" as the active tab).
$infraSheet.Activate()
$infraSheet.Range("H3").Select()

# Make "o&m" the active sheet (this also clears tabSelected on the
# previously-active "infrastructure" sheet).
$omSheet.Activate()

# Highlight row 26 with a solid red fill (RGB C00000), touching only the
# columns that already carry data so no new blank cells are materialized
# for the empty E/F columns.
$omSheet.Range("A26:D26").Interior.Color = 192
$omSheet.Range("G26:Q26").Interior.Color = 192

# Select the highlighted row on "o&m" (whole-row selection, active cell A26).
$omSheet.Range("A26:XFD26").Select()
